$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 with the new deal parameters.
$ws.Range("B1").Value = "AUD"
$ws.Range("E1").Value = "MasterCard"

# C1/D1/F1 hold numeric-looking text (amounts / card numbers) that must stay
# stored as text, not be auto-converted to numbers. Force text entry via a
# temporary "@" number format, then clear the format again so the cell keeps
# its default (unstyled) look, matching how the rest of the sheet is stored.
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "5098.97"
$ws.Range("C1").ClearFormats()

$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Value = "160"
$ws.Range("D1").ClearFormats()

$ws.Range("F1").NumberFormat = "@"
$ws.Range("F1").Value = "3388028672"
$ws.Range("F1").ClearFormats()

# Remove the second row entirely.
$ws.Rows.Item(2).Delete()
